$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EPBDS-12566: Numbers.toString() now renders +/-Infinity with the math
# symbol "\u221E" instead of the word "Infinity" / "-Infinity". Update the
# two "expected result" demo rows that show the toString() output for an
# infinite input: the miscDouble block (row 12) and the miscFloat block
# (row 24). C12/C24 ("= toString(v)") already carry the quote-prefix style
# used throughout these rows, so copy formats from there after writing the
# new text to keep each cell's style index intact.
$infinity = [char]0x221E
$negInfinity = "-" + [char]0x221E

$targets = "L12", "M12", "L24", "M24"
$values = $infinity, $negInfinity, $infinity, $negInfinity

for ($i = 0; $i -lt $targets.Length; $i++) {
    $cell = $ws.Range($targets[$i])
    $cell.Value = $values[$i]
    $ws.Range("C12").Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Selection moved off the stale T23 cell onto M23.
$ws.Range("M23").Select()
